$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.789.32'
$ws.Range('E2').Value = '  -0.47%  '

$ws.Range('D3').Value = '1.871.85'
$ws.Range('E3').Value = '  -0.25%  '

$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  +0.22%  '

$ws.Range('D5').Value = '''0.7187'
$ws.Range('E5').Value = '  -3.23%  '

$ws.Range('D6').Value = '''241.73'
$ws.Range('E6').Value = '  -0.42%  '

$ws.Range('E7').Value = '  +0.25%  '

$ws.Range('D8').Value = '''0.3140'
$ws.Range('E8').Value = '  -0.68%  '

$ws.Range('D9').Value = '''0.07446'
$ws.Range('E9').Value = '  +3.60%  '

$ws.Range('D10').Value = '''24.51'
$ws.Range('E10').Value = '  -1.26%  '

$ws.Range('D11').Value = '''0.08179'

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.892.80'
$ws.Range('E12').Value = '  +0.68%  '

$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = '''0.7427'
$ws.Range('E13').Value = '  -1.20%  '

$ws.Range('D14').Value = '''5.316'
$ws.Range('E14').Value = '  -2.45%  '

$ws.Range('D15').Value = '''92.42'
$ws.Range('E15').Value = '  -0.36%  '

$ws.Range('D16').Value = '29.903.21'
$ws.Range('E16').Value = '  -0.06%  '

$ws.Range('D17').Value = '''5.993'
$ws.Range('E17').Value = '  -1.83%  '

$ws.Range('D18').Value = '''246.20'
$ws.Range('E18').Value = '  +0.51%  '

$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '''0.000007905'
$ws.Range('E19').Value = '  +1.01%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '''13.45'
$ws.Range('E20').Value = '  -1.04%  '

$ws.Range('D21').Value = '2.152.17'
$ws.Range('E21').Value = '  +1.02%  '

$ws.Range('D22').Value = '''1.003'
$ws.Range('E22').Value = '  +0.37%  '

$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  +0.38%  '

$ws.Range('D24').Value = '''7.694'
$ws.Range('E24').Value = '  -3.99%  '

$ws.Range('D25').Value = '''9.188'
$ws.Range('E25').Value = '  -1.05%  '

$ws.Range('E26').Value = '  -3.96%  '

$ws.Range('D27').Value = '''163.48'
$ws.Range('E27').Value = '  -1.06%  '

$ws.Range('D28').Value = '''18.55'
$ws.Range('E28').Value = '  -0.47%  '

$ws.Range('D29').Value = '''1.996'
$ws.Range('E29').Value = '  -2.15%  '

$ws.Range('E30').Value = '  -6.88%  '

$ws.Range('E31').Value = '  -1.70%  '

$ws.Range('D32').Value = '''1.526'
$ws.Range('E32').Value = '  -0.65%  '

$ws.Range('D33').Value = '''4.161'
$ws.Range('E33').Value = '  -2.73%  '

$ws.Range('D34').Value = '''0.05390'
$ws.Range('E34').Value = '  +1.11%  '

$ws.Range('D35').Value = '''1.225'
$ws.Range('E35').Value = '  -1.18%  '

$ws.Range('D36').Value = '''0.7327'
$ws.Range('E36').Value = '  -2.91%  '

$ws.Range('D37').Value = '''0.9974'
$ws.Range('E37').Value = '  -0.30%  '

$ws.Range('D38').Value = '''2.697'
$ws.Range('E38').Value = '  -0.03%  '

$ws.Range('D39').Value = '''0.01911'
$ws.Range('E39').Value = '  -2.74%  '

$ws.Range('D40').Value = '''2.726'
$ws.Range('E40').Value = '  -0.96%  '

$ws.Range('D41').Value = '''0.4435'
$ws.Range('E41').Value = '  -2.25%  '

$ws.Range('D42').Value = '''0.8881'
$ws.Range('E42').Value = '  +3.67%  '

$ws.Range('D43').Value = '''5.977'
$ws.Range('E43').Value = '  -1.43%  '

$ws.Range('D44').Value = '''71.45'
$ws.Range('E44').Value = '  -1.73%  '

$ws.Range('E45').Value = '  +0.12%  '

$ws.Range('D46').Value = '1.039.87'
$ws.Range('E46').Value = '  -6.76%  '

$ws.Range('D47').Value = '''103.55'
$ws.Range('E47').Value = '  +0.00%  '

$ws.Range('D48').Value = '''7.447'
$ws.Range('E48').Value = '  -2.72%  '

$ws.Range('D49').Value = '''9.628'
$ws.Range('E49').Value = '  +0.39%  '

$ws.Range('D50').Value = '''1.798'
$ws.Range('E50').Value = '  -2.52%  '

$ws.Range('D51').Value = '2.037.97'
$ws.Range('E51').Value = '  +0.83%  '
